# The "Förändrad" (Changed) date column (C) was bumped by one day
# (2023-10-06 -> 2023-10-07, serial 45205 -> 45206) for every data row
# (rows 2 through 319) on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 319 }

$rng = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3))
$rng.Value = 45206
